$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = 44040
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "220907"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "18878"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = 34141
$ws.Range("F4").Value = 5263
$ws.Range("K4").Value = 114142
$ws.Range("L4").Value = 17290

# Row 8
$ws.Range("B8").Value = 44040
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "38855"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "286"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "953"
$ws.Range("E8").Style = "Normal"
$ws.Range("H8").Value = 1.75

# Row 10
$ws.Range("C10").Value = 40181
$ws.Range("D10").Value = 428
$ws.Range("E10").Value = 8424
$ws.Range("F10").Value = 110
$ws.Range("G10").Value = 24.35
$ws.Range("H10").Value = 25.88
$ws.Range("K10").Value = 34593
$ws.Range("L10").Value = 425

# Row 13
$ws.Range("B13").Value = 44040
$ws.Range("C13").Value = 19791
$ws.Range("D13").Value = 626
$ws.Range("E13").Value = 372

# Row 17
$ws.Range("B17").Value = 44039
$ws.Range("C17").Value = 178642
$ws.Range("D17").Value = 4426
$ws.Range("E17").Value = 4774
$ws.Range("F17").Value = 441
$ws.Range("G17").Value = 4.62
$ws.Range("H17").Value = 10.67
$ws.Range("K17").Value = 103302
$ws.Range("L17").Value = 4133

# Row 19
$ws.Range("B19").Value = 44039
$ws.Range("C19").Value = 54299
$ws.Range("D19").Value = 1543
$ws.Range("E19").Value = 23585
$ws.Range("F19").Value = 772

# Row 21
$ws.Range("B21").Value = 44040
$ws.Range("C21").Value = 106331
$ws.Range("D21").Value = 7146
$ws.Range("E21").Value = 15052
$ws.Range("G21").Value = 29.82
$ws.Range("K21").Value = 50476

# Row 23
$ws.Range("B23").Value = 44040
$ws.Range("C23").Value = 3475
$ws.Range("D23").Value = 51
$ws.Range("G23").Value = 0.58

# Row 25
$ws.Range("B25").Value = 44039
$ws.Range("C25").Value = 47089
$ws.Range("D25").Value = 3541
$ws.Range("E25").Value = 6443
$ws.Range("F25").Value = 654
$ws.Range("G25").Value = 13.68

# Row 27
$ws.Range("B27").Value = 44040
$ws.Range("C27").Value = 45314
$ws.Range("D27").Value = 1807
$ws.Range("E27").Value = 2130
$ws.Range("F27").Value = 121
$ws.Range("G27").Value = 6.07
$ws.Range("H27").Value = 6.93
$ws.Range("K27").Value = 35069
$ws.Range("L27").Value = 1745

# Row 28
$ws.Range("B28").Value = 44040
$ws.Range("C28").Value = 25157
$ws.Range("D28").Value = 321
$ws.Range("E28").Value = 1540
$ws.Range("G28").Value = 7.85
$ws.Range("H28").Value = 7.87
$ws.Range("K28").Value = 19625
$ws.Range("L28").Value = 305

# Row 29
$ws.Range("B29").Value = 44040
$ws.Range("C29").Value = 79090
$ws.Range("D29").Value = 6091
$ws.Range("E29").Value = 22020
$ws.Range("F29").Value = 2426
$ws.Range("G29").Value = 27.84
$ws.Range("H29").Value = 39.83

# Row 30
$ws.Range("C30").Value = 111038
$ws.Range("D30").Value = 3700

# Row 31
$ws.Range("B31").Value = 44039
$ws.Range("C31").Value = 466550
$ws.Range("D31").Value = 8518
$ws.Range("E31").Value = 12866
$ws.Range("G31").Value = 4.3
$ws.Range("H31").Value = 8.49
$ws.Range("K31").Value = 298953
$ws.Range("L31").Value = 8261

# Row 32
$ws.Range("B32").Value = 44040
$ws.Range("C32").Value = 63678
$ws.Range("D32").Value = 2725
$ws.Range("E32").Value = 7305
$ws.Range("G32").Value = 11.47
$ws.Range("H32").Value = 14.06

# Row 33
$ws.Range("B33").Value = 44040
$ws.Range("C33").Value = 1699
$ws.Range("D33").Value = 22
$ws.Range("E33").Value = 77
$ws.Range("G33").Value = 4.94
$ws.Range("K33").Value = 1560
$ws.Range("L33").Value = 22

# Row 34
$ws.Range("B34").Value = 44040
$ws.Range("C34").Value = 50179
$ws.Range("D34").Value = 906
$ws.Range("E34").Value = 7716
$ws.Range("F34").Value = 204
$ws.Range("G34").Value = 17.01
$ws.Range("H34").Value = 22.97
$ws.Range("K34").Value = 45352
$ws.Range("L34").Value = 888

# Row 35
$ws.Range("B35").Value = 44040
$ws.Range("C35").Value = 175052
$ws.Range("D35").Value = 3563
$ws.Range("E35").Value = 44977
$ws.Range("F35").Value = 1617
$ws.Range("G35").Value = 25.69
$ws.Range("H35").Value = 45.38

# Row 37
$ws.Range("B37").Value = 44040
$ws.Range("C37").Value = 6500
$ws.Range("E37").Value = 337
$ws.Range("K37").Value = 5630

# Row 41
$ws.Range("C41").Value = 42789
$ws.Range("D41").Value = 839
$ws.Range("E41").Value = 3404
$ws.Range("G41").Value = 7.96

# Row 42
$ws.Range("B42").Value = 44040
$ws.Range("C42").Value = 116087
$ws.Range("D42").Value = 1820
$ws.Range("E42").Value = 19790
$ws.Range("F42").Value = 569
$ws.Range("G42").Value = 24.04
$ws.Range("H42").Value = 32.38
$ws.Range("K42").Value = 82329
$ws.Range("L42").Value = 1757

# Row 43
$ws.Range("B43").Value = 44040
$ws.Range("C43").Value = 173731
$ws.Range("D43").Value = 7446
$ws.Range("E43").Value = 28884
$ws.Range("F43").Value = 2043
$ws.Range("H43").Value = 27.44

# Row 44
$ws.Range("B44").Value = 44040
$ws.Range("C44").Value = 19222
$ws.Range("D44").Value = 158
$ws.Range("E44").Value = 228
$ws.Range("G44").Value = 1.19
$ws.Range("H44").Value = 1.9

# Row 45
$ws.Range("B45").NumberFormat = "YYYY-MM-DD"
$ws.Range("B45").Value = 44040
$ws.Range("C45").Value = 86497
$ws.Range("D45").Value = 3382
$ws.Range("E45").Value = 21925
$ws.Range("F45").Value = 652
$ws.Range("G45").Value = 29.34
$ws.Range("H45").Value = 19.63
$ws.Range("J45").Value = $true
$ws.Range("K45").Value = 74731
$ws.Range("L45").Value = 3322
$ws.Range("O45").Value = "Success!"

# Row 46
$ws.Range("B46").Value = 44040
$ws.Range("C46").Value = 52281
$ws.Range("D46").Value = 1580
$ws.Range("E46").Value = 10776
$ws.Range("G46").Value = 20.61
$ws.Range("H46").Value = 9.43

# Row 47
$ws.Range("B47").Value = 44040
$ws.Range("C47").Value = 116182
$ws.Range("D47").Value = 8551
$ws.Range("E47").Value = 10870
$ws.Range("F47").Value = 700

# Row 48
$ws.Range("B48").Value = 44040
$ws.Range("C48").Value = 18485
$ws.Range("D48").Value = 381
$ws.Range("E48").Value = 5198
$ws.Range("G48").Value = 30.29
$ws.Range("H48").Value = 38.85
$ws.Range("K48").Value = 17161
$ws.Range("L48").Value = 381

# Row 49
$ws.Range("B49").Value = 44040
$ws.Range("C49").Value = 44819
$ws.Range("D49").Value = 1213
$ws.Range("E49").Value = 9479
$ws.Range("F49").Value = 402
$ws.Range("G49").Value = 27.96
$ws.Range("H49").Value = 34.39
$ws.Range("K49").Value = 33899
$ws.Range("L49").Value = 1169

# Row 50
$ws.Range("B50").Value = 44039
$ws.Range("C50").Value = 412878
$ws.Range("D50").Value = 25126
$ws.Range("F50").Value = 6356
$ws.Range("L50").Value = 23623

# Row 51
$ws.Range("B51").Value = 44040
$ws.Range("C51").Value = 84109
$ws.Range("D51").Value = 1565
$ws.Range("E51").Value = 28332
$ws.Range("F51").Value = 642
$ws.Range("G51").Value = 38.21
$ws.Range("H51").Value = 43.41
$ws.Range("K51").Value = 74148
$ws.Range("L51").Value = 1479
